$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Correct the manufacturer part number typo (LSH -> LTH)
$ws.Range("D2").Value = "LTH-030-01-X-D-A-TR"

# Leave selection on G8, matching the saved cursor position
$ws.Range("G8").Select()
